# Apply the "May 9th" data changes to the falling-data worksheet.
#
# What actually changed (columns C:H only - the sensor readings ax,ay,az,gx,gy,gz):
#   - 3 brand-new sensor readings are inserted right after the header row;
#     the sensor readings that used to occupy rows 2-21 are pushed down by
#     3 rows, ending up in rows 5-24 (same values, just moved down).
#   - 7 brand-new sensor readings are appended after the (now shifted) last
#     row of data, landing in rows 25-31.
#   - Column A (timestamp) and column B (label) are NOT touched by the
#     shift: they simply keep following the existing "timestamp = (row-2)*100,
#     label = 'falling'" pattern for every row, including the freshly
#     appended ones.
#   - The worksheet dimension grows from A1:H21 to A1:H31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: capture the current sensor readings (C2:H21) before we move them ---
# (use Value2 for the read - the Value getter in this runtime does not return
#  a usable array; Value2 does and is otherwise equivalent for numbers/strings)
$existing = $ws.Range("C2:H21").Value2

# --- Step 2: shift those sensor readings down by 3 rows -> C5:H24 ---
$ws.Range("C5:H24").Value2 = $existing

# --- Step 3: write the 3 brand-new sensor readings into the freed rows C2:H4 ---
$topRows = @(
    @(-0.01334476470947288, 0.1471533775329589,  -0.2691573500633239, -0.0074830991216003, 0.0424551330506801,  0.0247400421649217),
    @(-0.2625431418418885,  0.13736093044281,     -0.1842701695859431,  0.0018325957935303,  0.0197004042565822,  0.0239764600992202),
    @(-0.3044750690460205,  0.2236802577972413,   -0.08151795715093602, 0.0056505035609006, -0.0122173046693205,  0.009315694682300001)
)

$r = 2
foreach ($row in $topRows) {
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $r++
}

# --- Step 3b: the shift above only moved columns C:H, so rows 22-24 (which
#     used to be beyond the original A1:H21 data) still need their
#     timestamp/label filled in, continuing the same pattern used by every
#     other row ( timestamp = (row-2)*100, label = "falling" ) ---
for ($r = 22; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
    $ws.Cells.Item($r, 2).Value = 'falling'
}

# --- Step 4: append 7 brand-new rows (timestamp/label/sensor readings) after row 24 ---
$bottomRows = @(
    @(2300, 'falling',  0.09789943695068341,  0.3452561050653447,  -0.131537172943354,   0.0568104684352874,  0.1012509167194366,  -0.0526871271431446),
    @(2400, 'falling',  0.05506801605224545,  0.07766664028167738, -0.2464380264282224,   0.0032070425804704, -0.0478002056479454,  -0.024892758578062),
    @(2500, 'falling', -0.06795549392700258,  0.1180151626467706,  -0.2082828953862188,  -0.0117591563612222,  0.0113010071218013,   0.0294742472469806),
    @(2600, 'falling', -0.01404476165771439,  0.2834141030907641,  -0.1360972765833135,  -0.0128281703218817, -0.0500909499824047,  -0.0126754539087414),
    @(2700, 'falling', -0.01649236679077155,  0.2205449156463143,  -0.119759158231318,    0.0056505035609006, -0.0193949714303016,  -0.0198531206697225),
    @(2800, 'falling', -0.003359794616699139, 0.1425043791532516,  -0.1422623544931412,   0.0088575463742017,  0.0649044290184974,   0.0067195175215601),
    @(2900, 'falling',  0.02320241928100578,  0.1721755955368283,  -0.1290906090289352,  -0.0021380283869802,  0.01328631862998,     0.0001527163112768)
)

$r = 25
foreach ($row in $bottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
